$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 0.552020797677646
$ws.Range("C5").Value = 0.03721868832744345
$ws.Range("D5").Value = 0.04177234892042988
$ws.Range("B7").Value = 5.887746458637575
$ws.Range("C7").Value = 0.14033033230028319
$ws.Range("D7").Value = 1.5747381783615222
$ws.Range("B8").Value = 0.56427698494402989
$ws.Range("C8").Value = 0.04685886730259514
$ws.Range("D8").Value = 0.11816619899765356
$ws.Range("B11").Value = 1.6688105044077446
$ws.Range("C11").Value = 0.0313992299162109
$ws.Range("D11").Value = 0.24674824722578104
$ws.Range("B12").Value = 4.4556905606570183
$ws.Range("C12").Value = 0.11981751728257885
$ws.Range("D12").Value = 1.7179330027068458
$ws.Range("B14").Value = 3.0662983744488881
$ws.Range("C14").Value = 0.07829053470477146
$ws.Range("D14").Value = 0.71605070902052759
$ws.Range("B17").Value = 0.49300033522356285
$ws.Range("C17").Value = 0.0029742871745181377
$ws.Range("D17").Value = 0.1591624080087729
$ws.Range("B18").Value = 1.0595139569153285
$ws.Range("C18").Value = 0.03860165014442016
$ws.Range("D18").Value = 1.0132129735914934
$ws.Range("B19").Value = 1.2235329134554156
$ws.Range("C19").Value = 0.0830334860277174
$ws.Range("D19").Value = 0.88894514207060304
$ws.Range("B20").Value = 10.118623953953241
$ws.Range("C20").Value = 0.4972681892915769
$ws.Range("D20").Value = 2.0073197206277169
$ws.Range("B21").Value = 2.2015932773276417
$ws.Range("C21").Value = 0.10428905128668986
$ws.Range("D21").Value = 0.18871690821228074
$ws.Range("B22").Value = 1.1078115039356924
$ws.Range("C22").Value = 0.02140911269705824
$ws.Range("D22").Value = 0.45234337288385151
$ws.Range("B23").Value = 2.3256318403863818
$ws.Range("C23").Value = 0.09821911044352538
$ws.Range("D23").Value = 1.2638168626566357
$ws.Range("B24").Value = 0.42136009761146614
$ws.Range("C24").Value = 0.03316249852242005
$ws.Range("D24").Value = 0.052589511929746736
$ws.Range("B25").Value = 2.3964294277877904
$ws.Range("C25").Value = 0.02003887252744873
$ws.Range("D25").Value = 1.5524555838533873
$ws.Range("B26").Value = 1.91712384551855
$ws.Range("C26").Value = 0.10609531991554882
$ws.Range("D26").Value = 0.49088698085254556
$ws.Range("B27").Value = 3.2422529429629865
$ws.Range("C27").Value = 0.16022280778190898
$ws.Range("D27").Value = 0.73629706566311548
$ws.Range("B28").Value = 2.1417347584145414
$ws.Range("C28").Value = 0.09589363531582511
$ws.Range("D28").Value = 0.32655562559755169
